$wb = $excel.ActiveWorkbook

# "Generate Report for Handoff" - update the Latest Handoff Datetime (and the
# Overview sheet's Latest HO Xliff Generate Date) for the file that was just
# handed off: 42d7c95d-cfa1-4135-9c2b-b1731acd57b3.md

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$overview.Range("G4").Value = "2016-10-18 03:43:59"
$zhcn.Range("H4").Value = "2016-10-18 03:43:44"
$dede.Range("H4").Value = "2016-10-18 03:43:59"
